$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats constant used with PasteSpecial to copy only the cell's formatting
# (fill/font/alignment -> style index), leaving its value/formula untouched.
$xlPasteFormats = -4122

function Move-CardContent($WsParam, $SourceCell, $DestCell) {
    $src = $WsParam.Range($SourceCell)
    $dst = $WsParam.Range($DestCell)

    # Copy the source's formatting onto the destination cell.
    $src.Copy()
    $dst.PasteSpecial($xlPasteFormats)

    # Move the text value itself from the source cell to the destination cell.
    $dst.Value2 = $src.Value2
    $src.Value2 = ""
}

# Tasks progressing from "DOING" (column H) to "DONE" (column J).
Move-CardContent $ws "H10" "J10"
Move-CardContent $ws "H13" "J13"
Move-CardContent $ws "H16" "J16"

# Keep the "DOING" column's empty placeholder cell in sync with the row's text-style banding.
$ws.Range("H19").Copy()
$ws.Range("J19").PasteSpecial($xlPasteFormats)

# Tasks progressing from the backlog "STORY" column (D) into "DOING" (column H).
Move-CardContent $ws "D22" "H22"
Move-CardContent $ws "D25" "H25"

$excel.CutCopyMode = 0
